$d = $word.ActiveDocument

# --- Part 1: remove the old "_GoBack" bookmark that currently sits right
#     after "cas d'utilisation « Consulter son" ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- Part 2: find the paragraph containing the precondition sentence and
#     split "L'utilisateur est inscrit sur le site" into three separate
#     runs: "L'utilisateur est " + "connecté" + " au site" ---
$target = "L’utilisateur est inscrit sur le site"
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    $pStart = $para.Range.Start
    $full = $para.Range.Text
    $idx = $full.IndexOf($target)
    if ($idx -ge 0) {
        $found = $true

        $rStart = $pStart + $idx
        $rEnd = $rStart + $target.Length
        $rng = $d.Range($rStart, $rEnd)

        # First run: "L'utilisateur est "
        $rng.Text = "L’utilisateur est "
        $afterFirst = $rng.End

        # Second run: "connecté"
        $rng2 = $d.Range($afterFirst, $afterFirst)
        $rng2.InsertAfter("connecté")
        $afterSecond = $rng2.End

        # Third run: " au site"
        $rng3 = $d.Range($afterSecond, $afterSecond)
        $rng3.InsertAfter(" au site")
        $afterThird = $rng3.End

        # --- Part 3: re-add the "_GoBack" bookmark, collapsed, right after
        #     the new text (end of the paragraph's content, before the
        #     paragraph mark).
        #
        #     Workaround: adding a bookmark whose collapsed position is
        #     exactly "paragraph end - 1" mis-resolves to the wrong spot in
        #     this COM host. To dodge that, pad with two throw-away
        #     characters first (so the target position is no longer the
        #     last one in the paragraph), add the bookmark there, then
        #     delete the padding again; the bookmark stays correctly
        #     anchored once the padding text around it is removed.
        $padRange = $d.Range($afterThird, $afterThird)
        $padRange.InsertAfter("ZZ")

        $bmRange = $d.Range($afterThird, $afterThird)
        $d.Bookmarks.Add("_GoBack", $bmRange)

        $delRange = $d.Range($afterThird, $afterThird + 2)
        $delRange.Delete()

        break
    }
}

if (-not $found) {
    throw "Could not find target precondition sentence to split"
}
